# Update belief_mean values with final computed data (run prepare & render with final data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 40.5493307668479
$ws.Range("K2").Value = 36.9055302021312
$ws.Range("L2").Value = 36.9057907051389
$ws.Range("N2").Value = 41.6894250824717

$ws.Range("B3").Value = 34.7872415482579
$ws.Range("K3").Value = 30.4357942185907
$ws.Range("N3").Value = 43.0029533260978
